# Generate Report for Handback
# - Overview/zh-cn/de-de "Status"-type cells move from "Ready for handoff" to
#   "Handed back: in sync with en-US"
# - zh-cn + de-de rows now carry a "Latest Target File" (a.md, hyperlinked)
#   and a "Latest Handback File" / "Latest Handback DateTime" for file a.md
#   (each locale got its own handback datetime)

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28a3f882d89c73d7083775a34a5c25f110aefdb6/e2e/a.md"
$hyperlinkColor = 15570276  # BGR int for RGB FF6495ED (the workbook's existing HyperLink font color)

# ---------------------------------------------------------------------------
# Overview sheet: both rows' zh-cn / de-de status cells
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------------
# Helper: fill in one locale sheet (zh-cn or de-de) with its handback info
# ---------------------------------------------------------------------------
function Set-HandbackInfo($sheetName, $xlfName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column now reads the handed-back message for both data rows
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File -> a.md, with a hyperlink back to the source file
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = $hyperlinkColor

    # Latest Handback File + Latest Handback DateTime
    $ws.Range("J2").Value = $xlfName
    $ws.Range("J3").Value = $xlfName
    $ws.Range("K2").Value = $handbackDateTime
    $ws.Range("K3").Value = $handbackDateTime

    $ws.Columns.Item(3).ColumnWidth = 29.144371396019366
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Set-HandbackInfo "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-08-27 04:35:59"
Set-HandbackInfo "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-08-27 04:36:09"
